$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.522.51"
$ws.Range("E2").Value = "  +0.78%  "
$ws.Range("D3").Value = "1.879.04"
$ws.Range("E3").Value = "  +1.08%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'0.7176"
$ws.Range("E5").Value = "  +2.28%  "
$ws.Range("D6").Value = "'241.86"
$ws.Range("E6").Value = "  +1.57%  "
$ws.Range("D7").Value = "'1.002"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "'0.07940"
$ws.Range("E8").Value = "  +0.74%  "
$ws.Range("E9").Value = "  +2.32%  "
$ws.Range("D10").Value = "'25.38"
$ws.Range("E10").Value = "  +3.62%  "
$ws.Range("D11").Value = "'0.08275"
$ws.Range("D12").Value = "1.906.62"
$ws.Range("E12").Value = "  +1.99%  "
$ws.Range("D13").Value = "'0.7301"
$ws.Range("E13").Value = "  +3.31%  "
$ws.Range("D14").Value = "'5.287"
$ws.Range("E14").Value = "  +1.42%  "
$ws.Range("D15").Value = "'91.21"
$ws.Range("E15").Value = "  +1.88%  "
$ws.Range("D16").Value = "29.545.72"
$ws.Range("D17").Value = "'5.912"
$ws.Range("E17").Value = "  +1.74%  "
$ws.Range("D18").Value = "'246.22"
$ws.Range("E18").Value = "  +3.78%  "
$ws.Range("D19").Value = "'0.000007871"
$ws.Range("E19").Value = "  +0.57%  "
$ws.Range("D20").Value = "'13.33"
$ws.Range("E20").Value = "  +0.96%  "
$ws.Range("D21").Value = "2.120.88"
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("D22").Value = "'8.061"
$ws.Range("E22").Value = "  +6.44%  "
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").Value = "'0.1615"
$ws.Range("E25").Value = "  +13.94%  "
$ws.Range("D26").Value = "'163.66"
$ws.Range("E26").Value = "  +0.72%  "
$ws.Range("D27").Value = "'9.052"
$ws.Range("E27").Value = "  +1.72%  "
$ws.Range("D28").Value = "'18.33"
$ws.Range("E28").Value = "  +1.40%  "
$ws.Range("D29").Value = "'1.355"
$ws.Range("E29").Value = "  -3.12%  "
$ws.Range("E30").Value = "  +0.99%  "
$ws.Range("D31").Value = "'4.401"
$ws.Range("E31").Value = "  +2.37%  "
$ws.Range("D32").Value = "'4.107"
$ws.Range("E32").Value = "  +1.91%  "
$ws.Range("D33").Value = "'0.05217"
$ws.Range("E33").Value = "  +0.83%  "
$ws.Range("D34").Value = "'1.948"
$ws.Range("E34").Value = "  +2.24%  "
$ws.Range("D35").Value = "'1.199"
$ws.Range("E35").Value = "  +1.74%  "
$ws.Range("D36").Value = "'0.7275"
$ws.Range("E36").Value = "  +2.69%  "
$ws.Range("D37").Value = "'2.681"
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("D38").Value = "'0.01872"
$ws.Range("E38").Value = "  +1.22%  "
$ws.Range("D39").Value = "1.208.69"
$ws.Range("E39").Value = "  +6.00%  "
$ws.Range("D40").Value = "'2.709"
$ws.Range("E40").Value = "  +0.87%  "
$ws.Range("D41").Value = "'0.9103"
$ws.Range("E41").Value = "  -1.05%  "
$ws.Range("D42").Value = "'6.176"
$ws.Range("E42").Value = "  +3.72%  "
$ws.Range("D43").Value = "'73.76"
$ws.Range("E43").Value = "  +4.85%  "
$ws.Range("D44").Value = "'1.002"
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("D45").Value = "'102.45"
$ws.Range("E45").Value = "  -0.44%  "
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "2.021.29"
$ws.Range("E46").Value = "  -0.19%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "'0.5300"
$ws.Range("E47").Value = "  -0.32%  "
$ws.Range("E48").Value = "  +3.10%  "
$ws.Range("D49").Value = "'2.932"
$ws.Range("E49").Value = "  +9.67%  "
$ws.Range("D50").Value = "'9.310"
$ws.Range("E50").Value = "  +1.37%  "
$ws.Range("D51").Value = "'0.4327"

# Reset formatting on cells that needed a quote-prefix to stay text,
# so their style index matches the original (unstyled) cells.
$ws.Range("D4").ClearFormats()
$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D7").ClearFormats()
$ws.Range("D8").ClearFormats()
$ws.Range("D10").ClearFormats()
$ws.Range("D11").ClearFormats()
$ws.Range("D13").ClearFormats()
$ws.Range("D14").ClearFormats()
$ws.Range("D15").ClearFormats()
$ws.Range("D17").ClearFormats()
$ws.Range("D18").ClearFormats()
$ws.Range("D19").ClearFormats()
$ws.Range("D20").ClearFormats()
$ws.Range("D22").ClearFormats()
$ws.Range("D25").ClearFormats()
$ws.Range("D26").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("D28").ClearFormats()
$ws.Range("D29").ClearFormats()
$ws.Range("D31").ClearFormats()
$ws.Range("D32").ClearFormats()
$ws.Range("D33").ClearFormats()
$ws.Range("D34").ClearFormats()
$ws.Range("D35").ClearFormats()
$ws.Range("D36").ClearFormats()
$ws.Range("D37").ClearFormats()
$ws.Range("D38").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("D41").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("D43").ClearFormats()
$ws.Range("D44").ClearFormats()
$ws.Range("D45").ClearFormats()
$ws.Range("D47").ClearFormats()
$ws.Range("D49").ClearFormats()
$ws.Range("D50").ClearFormats()
$ws.Range("D51").ClearFormats()
